$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.795630931854248
$ws.Range("B1").Value = 6.226106643676758
$ws.Range("C1").Value = 5.425380706787109
$ws.Range("D1").Value = 6.305791854858398
$ws.Range("E1").Value = 3.713513612747192
